# Update self-assessment points ("P" column, F) for several competence rows.
# Rows 13, 14, 16, 20, 22, 23 go from 3 -> 4; row 26 goes from 2 -> 4.
# All downstream SUM/ROUND formulas (F34, F39, F44, C46, F46, ...) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F13").Value = 4
$ws.Range("F14").Value = 4
$ws.Range("F16").Value = 4
$ws.Range("F20").Value = 4
$ws.Range("F22").Value = 4
$ws.Range("F23").Value = 4
$ws.Range("F26").Value = 4

# Move the viewport / selection to reflect where the user ended up working.
$ws.Activate()
$excel.Goto($ws.Range("H43"), $true)
